$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Remove the 5 rows describing keyboard(<bool>), showCurve(...), showExtraCurve(...),
# showEvents(...), showBackgroundEvents(<bool>) - rows 99 through 103.
$ws.Range("A99:C103").EntireRow.Delete()

# Narrow column C slightly (27.02 -> 24.71 characters). The host's
# ColumnWidth setter snaps to its own internal character-width grid, so
# 24.0 is the input that lands closest to the 24.71 target.
$ws.Columns.Item(3).ColumnWidth = 24.0

# Move the active selection to A3 (was C101, now out of range after the delete).
$ws.Range("A3").Select()
